$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1) Update the slide-layout date placeholder text (04/11/2015 -> 24/11/2015)
$layout = $s.CustomLayout
$dateShape = $layout.Shapes.Item("Date Placeholder 1")
$dateShape.TextFrame.TextRange.Text = "24/11/2015"

# 2) "Rectángulo 377" (Nodo04, id 378): merge the " es un número " run and
#    the trailing "par" run into a single run reading " es un número par".
$parShape = $s.Shapes.Item("Rectángulo 377")
$parRange = $parShape.TextFrame.TextRange
$parText = $parRange.Text
$parMarkerStart = $parText.IndexOf(" es un número ") + 1
$parMarkerLength = $parText.Length - $parMarkerStart + 1
$parRange.Characters($parMarkerStart, $parMarkerLength).Text = " es un número par"

# 3) "Rectángulo 378" (Nodo04, id 379): merge the " es un número " run and
#    the trailing "impar" run into a single run reading " es un número impar".
$imparShape = $s.Shapes.Item("Rectángulo 378")
$imparRange = $imparShape.TextFrame.TextRange
$imparText = $imparRange.Text
$imparMarkerStart = $imparText.IndexOf(" es un número ") + 1
$imparMarkerLength = $imparText.Length - $imparMarkerStart + 1
$imparRange.Characters($imparMarkerStart, $imparMarkerLength).Text = " es un número impar"

# 4) Remove the now-unused "algunas características son" connector textbox
#    and its associated fourth-level node (and the two angled connectors
#    that linked to it) from the concept map.
$s.Shapes.Item("CuadroTexto 363").Delete()
$s.Shapes.Item("Rectángulo 364").Delete()
$s.Shapes.Item("Conector angular 350").Delete()
$s.Shapes.Item("Conector angular 351").Delete()
